# Refined metadata to be additional tab
#
# 1) Refresh the "time_taken" (column F) timestamps on the existing "data"
#    sheet to reflect the re-run panel query.
# 2) Add a new "metadata" sheet (placed right after "data") describing the
#    panel query itself: name / id / version / version-created date /
#    query time / GET request URL.

$wb = $excel.ActiveWorkbook
$ds = $wb.Worksheets.Item("data")

$timestamps = @(
  "2021-10-05 14:23:05.691047",
  "2021-10-05 14:23:05.691057",
  "2021-10-05 14:23:05.691060",
  "2021-10-05 14:23:05.691063",
  "2021-10-05 14:23:05.691066",
  "2021-10-05 14:23:05.691069",
  "2021-10-05 14:23:05.691072",
  "2021-10-05 14:23:05.691075",
  "2021-10-05 14:23:05.691078",
  "2021-10-05 14:23:05.691081",
  "2021-10-05 14:23:05.691083",
  "2021-10-05 14:23:05.691086",
  "2021-10-05 14:23:05.691089",
  "2021-10-05 14:23:05.691091",
  "2021-10-05 14:23:05.691094",
  "2021-10-05 14:23:05.691097",
  "2021-10-05 14:23:05.691100",
  "2021-10-05 14:23:05.691102",
  "2021-10-05 14:23:05.691105",
  "2021-10-05 14:23:05.691108",
  "2021-10-05 14:23:05.691111",
  "2021-10-05 14:23:05.691114",
  "2021-10-05 14:23:05.691116",
  "2021-10-05 14:23:05.691119",
  "2021-10-05 14:23:05.691122",
  "2021-10-05 14:23:05.691125",
  "2021-10-05 14:23:05.691128",
  "2021-10-05 14:23:05.691131",
  "2021-10-05 14:23:05.691133",
  "2021-10-05 14:23:05.691136",
  "2021-10-05 14:23:05.691139",
  "2021-10-05 14:23:05.691142",
  "2021-10-05 14:23:05.691145",
  "2021-10-05 14:23:05.691148",
  "2021-10-05 14:23:05.691150",
  "2021-10-05 14:23:05.691153",
  "2021-10-05 14:23:05.691156",
  "2021-10-05 14:23:05.691159",
  "2021-10-05 14:23:05.691162",
  "2021-10-05 14:23:05.691164",
  "2021-10-05 14:23:05.691168",
  "2021-10-05 14:23:05.691170",
  "2021-10-05 14:23:05.691173",
  "2021-10-05 14:23:05.691176",
  "2021-10-05 14:23:05.691178",
  "2021-10-05 14:23:05.691181",
  "2021-10-05 14:23:05.691184",
  "2021-10-05 14:23:05.691187",
  "2021-10-05 14:23:05.691189",
  "2021-10-05 14:23:05.691192",
  "2021-10-05 14:23:05.691195",
  "2021-10-05 14:23:05.691197",
  "2021-10-05 14:23:05.691201",
  "2021-10-05 14:23:05.691203",
  "2021-10-05 14:23:05.691206",
  "2021-10-05 14:23:05.691209",
  "2021-10-05 14:23:05.691211",
  "2021-10-05 14:23:05.691214",
  "2021-10-05 14:23:05.691217",
  "2021-10-05 14:23:05.691219",
  "2021-10-05 14:23:05.691222"
)

$row = 2
foreach ($ts in $timestamps) {
    $ds.Cells.Item($row, 6).Value = $ts
    $row = $row + 1
}

# --- add the "metadata" sheet right after "data" ---
$md = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ds)
$md.Name = "metadata"

# match the bold/bordered/centered header look used on the "data" sheet
$ds.Range("B1").Copy()
$md.Range("B1:G1").PasteSpecial(-4122)

# match the index-column look used for A2 on the "data" sheet
$ds.Range("A2").Copy()
$md.Range("A2").PasteSpecial(-4122)

$md.Range("B1").Value = "data_name"
$md.Range("C1").Value = "data_id"
$md.Range("D1").Value = "data_version"
$md.Range("E1").Value = "data_version_created"
$md.Range("F1").Value = "panel_query_time"
$md.Range("G1").Value = "panel_get_request"

$md.Range("A2").Value = 0
$md.Range("B2").Value = "VACTERL-like phenotypes"
$md.Range("C2").Value = 101
$md.Range("D2").NumberFormat = "@"
$md.Range("D2").Value = "1.32"
$md.Range("E2").Value = "2021-07-08T10:47:25.122773Z"
$md.Range("F2").Value = "2021-10-05 14:23:05.687379"
$md.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/101/?format=json"

# keep "data" as the active/visible tab, as in the original workbook
$ds.Activate()

Write-Output "metadata sheet added; data!F2:F62 timestamps refreshed"
